# Commit: "api and order item"
#  1. Create Order: drop "order_id:string," from the order_items parameter list.
#  2. Update OrderItem: merge the "parameter:" / "{position:...}" runs (drop the space
#     between them).
#  3. Update OrderItem return value: merge "{" / "return:integer" / ",content:{...}}"
#     runs into one.
#  4. Two new (Word-generated) list-label character styles show up in styles.xml as a
#     side effect of the list-item edits above -- mint them explicitly.

$d = $word.ActiveDocument

# 1) Create Order -- parameter: {order:{user_id:string},order_items:[{order_id:string,position:string,part_id:string,quantity:string}]}
$d.Content.Find.Execute(
    "order_items:[{order_id:string,position:string,part_id:string,quantity:string}]}",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "order_items:[{position:string,part_id:string,quantity:string}]}",
    2) | Out-Null

# 2) Update OrderItem -- "parameter:" + "{position:string,part_id:string,quantity:string}"
#    collapse into a single run with no space in between.
$d.Content.Find.Execute(
    "parameter:{position:string,part_id:string,quantity:string}",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "parameter:{position:string,part_id:string,quantity:string}",
    2) | Out-Null

# 3) Update OrderItem return value -- "{" + "return:integer" + ",content:{...}}"
#    collapse into a single run.
$d.Content.Find.Execute(
    "{return:integer,content:{id:string,order_id:string,location_id:string,whouse_id:string,source_id:string,user_id:string,part_id:string,part_type:string,quantity:string}}",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "{return:integer,content:{id:string,order_id:string,location_id:string,whouse_id:string,source_id:string,user_id:string,part_id:string,part_type:string,quantity:string}}",
    2) | Out-Null

# 4) Mint the two new list-label character styles (ListLabel 63 / 64) that Word
#    generates for the edited numbered-list paragraphs.
$s63 = $d.Styles.Add("ListLabel63", 2)
$s63.NameLocal = "ListLabel 63"
$s63.Font.NameBi = "Symbol"

$s64 = $d.Styles.Add("ListLabel64", 2)
$s64.NameLocal = "ListLabel 64"
$s64.Font.NameBi = "OpenSymbol"
